$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

$newValues = @(
    "543_柔丽思红_undefined_undefined_1bunch",
    "738_鸢尾花蓝_undefined_undefined_1bunch",
    "739_鸢尾花白_undefined_undefined_1bunch",
    "400_龙船花红_undefined_undefined_1bunch",
    "705_中华桔梗紫_undefined_undefined_1bunch",
    "819_中华桔梗白_undefined_undefined_1bunch",
    "454_蓝星花_tweedia blue_undefined_1bunch",
    "476_酸浆果_Physalis_undefined_1bunch",
    "553_嘉兰百合`nmulti-colors_ gloriosa_undefined_1bunch",
    "447_黄金球_craspedia_undefined_1bunch"
)

$startRow = 22
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 3).Value = $newValues[$i]
}

$summary = $wb.Worksheets.Item("Summary")
$summary.Range("G2").NumberFormat = "@"
$summary.Cells.Item(2, 7).Value = "05200000000000000000000000000000"
